# truth-table.xlsx edit
#
# - "Sheet1" -> "TruthTable" (first tab, now the active tab)
# - a duplicate of it, "TruthTableGenerated", is inserted right after it
# - "Sheet2" (the original second sheet) stays last, untouched
# - TruthTable / TruthTableGenerated: B4 becomes the text "d",
#   B5 becomes the text "p" (was 0 / 1 respectively)

$wb = $excel.ActiveWorkbook

$truthTable = $wb.Worksheets.Item(1)

# Rename the first sheet.
$truthTable.Name = "TruthTable"

# Update the generated truth-table values (0/1 -> text labels).
$truthTable.Range("B4").Value = "d"
$truthTable.Range("B5").Value = "p"

# Duplicate TruthTable; the copy lands immediately after the source sheet.
$truthTable.Copy($null, $truthTable)
$generated = $wb.Worksheets.Item(2)
$generated.Name = "TruthTableGenerated"

# Leave the original Sheet2 as-is (now the 3rd / last tab).

# Make TruthTable the active sheet/tab again (Copy() activates the new copy).
$truthTable.Activate()
